# Auto commit at 2025-10-16  7:48:37.02
# Append two new daily rows (date 2025-10-15 / serial 45945) for both
# stations ("四方坪站充电量(kw)" and "高岭站充电量(kw)") to the bottom of
# Sheet1, following the existing layout/pattern used by the rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newDate = 45945

# Row 90 - 四方坪站充电量(kw)
$row90 = 90
$ws.Cells.Item($row90, 1).Value = $newDate
$ws.Cells.Item($row90, 2).Value = "四方坪站充电量(kw)"

$row90Values = @(
    575.55399999999997,
    1022.312,
    369.726,
    401.375,
    380.48500000000001,
    576.86099999999999,
    607.52699999999993,
    292.76,
    210.85400000000001,
    73.7,
    143.511,
    271.45699999999999,
    648.15999999999985,
    1354.2189999999996,
    524.53800000000001,
    245.34999999999997,
    199.75399999999999,
    165.27199999999999,
    50.06,
    88.57,
    105.58000000000001,
    94.049000000000007,
    0,
    138.28
)

for ($i = 0; $i -lt $row90Values.Length; $i++) {
    $ws.Cells.Item($row90, 3 + $i).Value = $row90Values[$i]
}

# Row 91 - 高岭站充电量(kw)
$row91 = 91
$ws.Cells.Item($row91, 1).Value = $newDate
$ws.Cells.Item($row91, 2).Value = "高岭站充电量(kw)"

$row91Values = @(
    327.19999999999993,
    392.48599999999993,
    91.914000000000016,
    153.24600000000001,
    27.390999999999998,
    46.585000000000001,
    481.11799999999994,
    206.77500000000003,
    242.19200000000001,
    173.62200000000001,
    30.800999999999998,
    265.14599999999996,
    292.88099999999997,
    484.23799999999989,
    419.62700000000001,
    245.23400000000001,
    181.45099999999999,
    112.736,
    0,
    71.274000000000001,
    7.0549999999999997,
    0,
    46.429000000000002,
    60.953000000000003
)

for ($i = 0; $i -lt $row91Values.Length; $i++) {
    $ws.Cells.Item($row91, 3 + $i).Value = $row91Values[$i]
}

# (Columns A and C:Z already carry a column-level style (date / 2-decimal
# number format) via <cols>, so newly written cells in those columns pick
# it up automatically - no explicit NumberFormat assignment needed, which
# also avoids minting duplicate numFmt/cellXfs entries.)

# Update the saved view state to match the scrolled-down position after
# the new rows were appended.
$ws.Application.ActiveWindow.ScrollRow = 62
$ws.Range("I94").Select()
